{"js": "// Office.js (Word JavaScript API) script.\n// Rewrites the body paragraphs following the \"1. presentation_input\" heading\n// to the revised create_summary_prompt copy (7 paragraphs instead of 5),\n// preserving the trailing manual line break on the final paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load('items');\nawait context.sync();\n\n// paragraphs.items[0]  = heading \"1. presentation_input\" (unchanged)\n// paragraphs.items[1..5] = the original 5 body paragraphs (jc=both)\n// Two more body paragraphs are needed (7 total), so append 2 new ones\n// after the last existing body paragraph before rewriting the text.\nconst lastBodyParagraph = paragraphs.items[5];\nlastBodyParagraph.insertParagraph('', Word.InsertLocation.after);\nlastBodyParagraph.insertParagraph('', Word.InsertLocation.after);\nawait context.sync();\n\nparagraphs.load('items');\nawait context.sync();\n\nparagraphs.items[1].insertText(`Artificial intelligence (AI) is rapidly transforming our world, and at the heart of this transformation lie large language models (LLMs).  AI, broadly defined, is a branch of computer science dedicated to creating intelligent systems capable of mimicking human cognitive functions. This involves a complex interplay of machine learning, deep learning, and natural language processing (NLP).  The applications of AI are vast and varied, ranging from the development of sophisticated robotics and autonomous vehicles to the creation of advanced decision-making systems used in various industries.`, Word.InsertLocation.replace);\nparagraphs.items[2].insertText(`Within the broader field of AI, LLMs represent a significant advancement. These are deep learning models trained on massive datasets of text and code.  Their architecture, typically based on the transformer model, allows them to process and generate human-like text with remarkable fluency and coherence.  Prominent examples of LLMs include GPT (Generative Pre-trained Transformer), BERT (Bidirectional Encoder Representations from Transformers), and LLaMA (Large Language Model Meta AI).`, Word.InsertLocation.replace);\nparagraphs.items[3].insertText(`The power of LLMs stems from a two-stage process: pretraining and fine-tuning.  Pretraining involves exposing the model to an enormous volume of text data, allowing it to learn the underlying patterns and relationships within language.  This initial training provides the model with a broad understanding of grammar, semantics, and even some aspects of world knowledge.  Fine-tuning then tailors the pre-trained model to specific tasks, such as question answering, text summarization, or machine translation.  This targeted training refines the model's performance on the desired application.`, Word.InsertLocation.replace);\nparagraphs.items[4].insertText(`A crucial aspect of LLM functionality is tokenization.  Before processing, the input text is broken down into smaller units, or tokens, which can be individual words, parts of words, or even sub-word units. This process allows the model to handle the complexities of language more effectively.  Furthermore, LLMs exhibit context awareness, meaning they can understand and generate text that is coherent and relevant to the preceding context. This ability is essential for creating natural-sounding conversations and generating meaningful responses.  The applications of LLMs are equally diverse, encompassing chatbots, machine translation services, content generation tools, and much more.`, Word.InsertLocation.replace);\nparagraphs.items[5].insertText(`However, the development and deployment of LLMs are not without challenges.  One significant concern is bias.  Because LLMs are trained on existing data, they can inherit and perpetuate biases present in that data, leading to unfair or discriminatory outputs.  Addressing this bias requires careful curation of training data and the development of techniques to mitigate biased outcomes.  Another challenge is the substantial computational cost associated with training and deploying LLMs.  The sheer scale of the data and the complexity of the models demand significant computing resources, making them expensive to develop and operate.`, Word.InsertLocation.replace);\nparagraphs.items[6].insertText(`Finally, ethical considerations are paramount.  The potential for misuse of LLMs, such as the generation of misinformation or the creation of deepfakes, raises serious ethical concerns.  Responsible development and deployment of LLMs require careful consideration of these ethical implications and the implementation of safeguards to prevent harmful applications.`, Word.InsertLocation.replace);\nawait context.sync();\n\n// The final paragraph also carries the trailing manual line break (<w:br/>)\n// that used to terminate the document's last paragraph, so it is inserted\n// via OOXML to keep the <w:t>/<w:br/> pair inside a single run, matching\n// the original authoring pattern.\nconst finalParagraph = paragraphs.items[7];\nconst finalRange = finalParagraph.getRange(Word.RangeLocation.content);\nconst finalOoxml = `<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:jc w:val=\"both\"/></w:pPr><w:r><w:t>In conclusion, AI and LLMs are transforming technology and society at an unprecedented pace.  Their potential applications are vast and transformative, but their development and deployment must be guided by a commitment to addressing the inherent challenges and ethical considerations.  Ongoing research is crucial to improve the efficiency, fairness, and safety of these powerful technologies, ensuring they are used for the benefit of humanity.</w:t><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>`;\nfinalRange.insertOoxml(finalOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Office Word COM interop script (PowerShell-style)\n# Rewrites the body paragraphs following the \"1. presentation_input\" heading\n# to the revised create_summary_prompt copy (7 paragraphs instead of 5),\n# preserving the trailing manual line break on the final paragraph.\n$d = $word.ActiveDocument\n\n# Paragraph 1 (Item 2) is unchanged in count; items 2-6 hold the original 5\n# body paragraphs. We need 7 body paragraphs total, so insert 2 new ones\n# after the last existing body paragraph (Item 6) before rewriting text.\n$lastBodyPara = $d.Paragraphs.Item(6)\n$lastBodyPara.Range.InsertParagraphAfter()\n$d.Paragraphs.Item(7).Range.InsertParagraphAfter()\n\n$d.Paragraphs.Item(2).Range.Text = \"Artificial intelligence (AI) is rapidly transforming our world, and at the heart of this transformation lie large language models (LLMs).  AI, broadly defined, is a branch of computer science dedicated to creating intelligent systems capable of mimicking human cognitive functions. This involves a complex interplay of machine learning, deep learning, and natural language processing (NLP).  The applications of AI are vast and varied, ranging from the development of sophisticated robotics and autonomous vehicles to the creation of advanced decision-making systems used in various industries.\"\n$d.Paragraphs.Item(3).Range.Text = \"Within the broader field of AI, LLMs represent a significant advancement. These are deep learning models trained on massive datasets of text and code.  Their architecture, typically based on the transformer model, allows them to process and generate human-like text with remarkable fluency and coherence.  Prominent examples of LLMs include GPT (Generative Pre-trained Transformer), BERT (Bidirectional Encoder Representations from Transformers), and LLaMA (Large Language Model Meta AI).\"\n$d.Paragraphs.Item(4).Range.Text = \"The power of LLMs stems from a two-stage process: pretraining and fine-tuning.  Pretraining involves exposing the model to an enormous volume of text data, allowing it to learn the underlying patterns and relationships within language.  This initial training provides the model with a broad understanding of grammar, semantics, and even some aspects of world knowledge.  Fine-tuning then tailors the pre-trained model to specific tasks, such as question answering, text summarization, or machine translation.  This targeted training refines the model's performance on the desired application.\"\n$d.Paragraphs.Item(5).Range.Text = \"A crucial aspect of LLM functionality is tokenization.  Before processing, the input text is broken down into smaller units, or tokens, which can be individual words, parts of words, or even sub-word units. This process allows the model to handle the complexities of language more effectively.  Furthermore, LLMs exhibit context awareness, meaning they can understand and generate text that is coherent and relevant to the preceding context. This ability is essential for creating natural-sounding conversations and generating meaningful responses.  The applications of LLMs are equally diverse, encompassing chatbots, machine translation services, content generation tools, and much more.\"\n$d.Paragraphs.Item(6).Range.Text = \"However, the development and deployment of LLMs are not without challenges.  One significant concern is bias.  Because LLMs are trained on existing data, they can inherit and perpetuate biases present in that data, leading to unfair or discriminatory outputs.  Addressing this bias requires careful curation of training data and the development of techniques to mitigate biased outcomes.  Another challenge is the substantial computational cost associated with training and deploying LLMs.  The sheer scale of the data and the complexity of the models demand significant computing resources, making them expensive to develop and operate.\"\n$d.Paragraphs.Item(7).Range.Text = \"Finally, ethical considerations are paramount.  The potential for misuse of LLMs, such as the generation of misinformation or the creation of deepfakes, raises serious ethical concerns.  Responsible development and deployment of LLMs require careful consideration of these ethical implications and the implementation of safeguards to prevent harmful applications.\"\n$d.Paragraphs.Item(8).Range.Text = \"In conclusion, AI and LLMs are transforming technology and society at an unprecedented pace.  Their potential applications are vast and transformative, but their development and deployment must be guided by a commitment to addressing the inherent challenges and ethical considerations.  Ongoing research is crucial to improve the efficiency, fairness, and safety of these powerful technologies, ensuring they are used for the benefit of humanity.\"\n\n# Restore the trailing manual line break (<w:br/>) that followed the final\n# paragraph in the original document, now on the new last paragraph.\n$lastPara = $d.Paragraphs.Item(8)\n$lastPara.Range.Find.Execute(\"of humanity.\", $false, $false, $false, $false, $false, $true, 1, $false, \"of humanity.^l\", 2) | Out-Null\n\n"}
